$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 81
$ws.Range("I2").Value = 185
$ws.Range("J2").Value = 801
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 199
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 140
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 4
$ws.Range("S2").Value = 88
$ws.Range("T2").Value = 119
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 1114
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1170
$ws.Range("AA2").Value = 9
